# Tratamento de materiais parte 1
# Applies updated quantity (G) and value (H) figures per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 417.000
$ws.Range("H2").Value = 1508.71
$ws.Range("G3").Value = 1866.000
$ws.Range("H3").Value = 244563.86
$ws.Range("G11").Value = 3984.000
$ws.Range("H11").Value = 122227.77
$ws.Range("G12").Value = 78.000
$ws.Range("H12").Value = 6660.45
$ws.Range("G16").Value = 1639.000
$ws.Range("H16").Value = 140081.29
$ws.Range("G20").Value = 2663.000
$ws.Range("H20").Value = 225401.66
$ws.Range("G23").Value = 67.000
$ws.Range("H23").Value = 9779.87
$ws.Range("G24").Value = 149.000
$ws.Range("H24").Value = 67351.70
$ws.Range("G28").Value = 469.480
$ws.Range("H28").Value = 17748.73
$ws.Range("G29").Value = 192.890
$ws.Range("H29").Value = 4944.38
$ws.Range("G30").Value = 620.600
$ws.Range("H30").Value = 7268.09
$ws.Range("G34").Value = 286.000
$ws.Range("H34").Value = 5306.73
$ws.Range("G45").Value = 10.000
$ws.Range("H45").Value = 975.81
$ws.Range("G49").Value = 22.000
$ws.Range("H49").Value = 2430.71
$ws.Range("G51").Value = 4.000
$ws.Range("H51").Value = 229.34
$ws.Range("G55").Value = 195.900
$ws.Range("H55").Value = 4208.31
$ws.Range("H60").Value = 3716.19
$ws.Range("H61").Value = 1707.25
$ws.Range("G62").Value = 14.000
$ws.Range("H62").Value = 7926.11
$ws.Range("G72").Value = 441.000
$ws.Range("H72").Value = 9934.13
$ws.Range("G76").Value = 16.000
$ws.Range("H76").Value = 965.23
$ws.Range("G77").Value = 2.000
$ws.Range("H77").Value = 46.69
$ws.Range("H79").Value = 11556.75
$ws.Range("H80").Value = 5495.51
$ws.Range("G87").Value = 4024.000
$ws.Range("H87").Value = 73399.08
$ws.Range("G95").Value = 27.000
$ws.Range("H95").Value = 410.66
$ws.Range("G96").Value = 116.000
$ws.Range("H96").Value = 662.62
$ws.Range("G102").Value = 46.000
$ws.Range("H102").Value = 6046.42
$ws.Range("H103").Value = 3803.61
$ws.Range("G107").Value = 429.000
$ws.Range("H107").Value = 1810.86
$ws.Range("G108").Value = 40.000
$ws.Range("H108").Value = 1106.48
$ws.Range("G109").Value = 44.000
$ws.Range("H109").Value = 1485.62
$ws.Range("H110").Value = 5565.89
$ws.Range("G111").Value = 177.000
$ws.Range("H111").Value = 6562.07
$ws.Range("G122").Value = 3.000
$ws.Range("H122").Value = 454.82
$ws.Range("G124").Value = 5.000
$ws.Range("H124").Value = 143.07
$ws.Range("G125").Value = 602.000
$ws.Range("H125").Value = 3387.57
$ws.Range("G128").Value = 7.000
$ws.Range("H128").Value = 1735.84
$ws.Range("G129").Value = 24.000
$ws.Range("H129").Value = 2489.07
$ws.Range("G133").Value = 27.000
$ws.Range("H133").Value = 5231.40
$ws.Range("H136").Value = 84289.17
$ws.Range("G137").Value = 174.000
$ws.Range("H137").Value = 27586.42
$ws.Range("G140").Value = 194.000
$ws.Range("H140").Value = 9113.09
$ws.Range("H173").Value = 6973.58
$ws.Range("G175").Value = 230.000
$ws.Range("H175").Value = 39693.26
$ws.Range("G176").Value = 262.000
$ws.Range("H176").Value = 8755.57
$ws.Range("G181").Value = 10.000
$ws.Range("H181").Value = 374.15
$ws.Range("H186").Value = 2079.32
$ws.Range("G187").Value = 455.000
$ws.Range("H187").Value = 3278.07
$ws.Range("G189").Value = 781.000
$ws.Range("H189").Value = 5262.42
$ws.Range("G192").Value = 3461.000
$ws.Range("H192").Value = 34794.15
$ws.Range("G200").Value = 23.000
$ws.Range("H200").Value = 11511.66
$ws.Range("G210").Value = 141.000
$ws.Range("H210").Value = 190.52
$ws.Range("G211").Value = 69575.000
$ws.Range("H211").Value = 163103.72
$ws.Range("G213").Value = 80.000
$ws.Range("H213").Value = 2851.51
$ws.Range("H218").Value = 1181.25
$ws.Range("H220").Value = 2387.60
$ws.Range("G221").Value = 288.000
$ws.Range("H221").Value = 15189.17
$ws.Range("G222").Value = 6082.000
$ws.Range("H222").Value = 51454.54
$ws.Range("G225").Value = 101.000
$ws.Range("H225").Value = 27673.49
$ws.Range("G229").Value = 194.000
$ws.Range("H229").Value = 13823.42
$ws.Range("G237").Value = 5.000
$ws.Range("H237").Value = 128.99
$ws.Range("G239").Value = 85.000
$ws.Range("H239").Value = 8045.60
$ws.Range("H240").Value = 5429.83
$ws.Range("G245").Value = 19.000
$ws.Range("H245").Value = 2491.66
$ws.Range("H255").Value = 2328355.94

